$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill blank 'Extracurricular Opportunities' (L) responses with 'Adequate'
$lRows = @(3,4,6,7,8,9,11,12,14,16,17,19,20,24,25,26,28,29,30,31,33,35,36,37,38,39,40,41,45,47,48,49,54,55,56,57,58,59,61,63,64,65,67,69,73,75,76,78,80,81,82,83,85,87,88,96,98)
foreach ($r in $lRows) { $ws.Range("L$r").Value = "Adequate" }

# Fill blank 'Campus Safety' (M) responses with 'Excellent'
$mRows = @(3,4,6,7,9,11,12,14,16,17,19,20,24,25,26,28,29,30,33,35,36,38,40,41,45,47,48,49,54,55,56,58,59,61,64,65,67,69,73,75,76,78,80,81,82,83,85,87,88,96)
foreach ($r in $mRows) { $ws.Range("M$r").Value = "Excellent" }

$wb.Save()
